$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header "pe" -> "pe_cl"
$ws.Range("B1").Value = "pe_cl"

# Add new column header "pe_jim"
$ws.Range("C1").Value = "pe_jim"

# Update column B values (new scale)
$ws.Range("B2").Value = 1
$ws.Range("B3").Value = 2
$ws.Range("B4").Value = 2

# Add column C values
$ws.Range("C2").Value = 1
$ws.Range("C3").Value = 2
$ws.Range("C4").Value = 2

# Update selection to match target state
$ws.Range("C2").Select()
